# Add a new Sprint 2 backlog item: "Remove console output"
# This corresponds to inserting a new row above the old (blank) row 28,
# which pushes the trailing rows (old 28-36) down to 29-37 and keeps all
# of the existing SUM()/difference formulas self-adjusting to the new
# ranges, exactly like using Excel's "Insert Sheet Rows" command.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sprint 2")
$ws.Activate()

# Insert a new blank row at position 28 (shifts 28..36 -> 29..37).
$ws.Rows("28:28").Insert()

# Row 27 has the exact same column formatting pattern (styles) that the
# new task row needs, so copy its formats down into the freshly inserted
# row 28 before filling in the new values.
$ws.Range("A27:G27").Copy()
$ws.Range("A28:G28").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# Fill in the new backlog item data.
$ws.Range("A28").Value = "1,2,3"
$ws.Range("B28").Value = 23
$ws.Range("C28").Value = "Remove console output"
$ws.Range("D28").Value = 0.25
$ws.Range("E28").Value = "AE"

# Match the saved view/selection state: Sprint 2 tab active, F28 selected.
$ws.Range("F28").Select()
